$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 22.70222272539856
$ws.Range("C2").Value = 3.819600466652046
$ws.Range("D2").Value = 4.589169096951085
$ws.Range("F2").Value = 54.86760364773193
$ws.Range("G2").Value = 3.791704092098842
$ws.Range("J2").Value = 10.83336959216128
$ws.Range("K2").Value = 18.48738854425894
$ws.Range("L2").Value = 11.8897750698502
$ws.Range("N2").Value = 26.02863754829579

$ws.Range("B3").Value = 22.59509705840143
$ws.Range("C3").Value = 3.646408588272827
$ws.Range("D3").Value = 4.596911052096146
$ws.Range("F3").Value = 54.82522718275818
$ws.Range("G3").Value = 3.794908304184778
$ws.Range("J3").Value = 10.85196087598211
$ws.Range("K3").Value = 18.41762434232059
$ws.Range("L3").Value = 11.902778269244
$ws.Range("N3").Value = 26.06494950888041

$ws.Range("B4").Value = 22.53460620208896
$ws.Range("C4").Value = 3.550785219245515
$ws.Range("D4").Value = 4.602302241292524
$ws.Range("F4").Value = 54.80855162299557
$ws.Range("G4").Value = 3.79697909113301
$ws.Range("J4").Value = 10.86422404617488
$ws.Range("K4").Value = 18.37901317437258
$ws.Range("L4").Value = 11.91254530345754
$ws.Range("N4").Value = 26.08910746240619

$ws.Range("B5").Value = 22.51130415336243
$ws.Range("C5").Value = 3.511387507067019
$ws.Range("D5").Value = 4.604659952741948
$ws.Range("F5").Value = 54.80410905964022
$ws.Range("G5").Value = 3.797849043125195
$ws.Range("J5").Value = 10.86943505862132
$ws.Range("K5").Value = 18.36435270810898
$ws.Range("L5").Value = 11.91697409641503
$ws.Range("N5").Value = 26.09942039758225

$ws.Range("B6").Value = 22.50751684442817
$ws.Range("C6").Value = 3.504821752875753
$ws.Range("D6").Value = 4.605061170687037
$ws.Range("F6").Value = 54.80351355395411
$ws.Range("G6").Value = 3.797995076430651
$ws.Range("J6").Value = 10.87031326091775
$ws.Range("K6").Value = 18.3619835241757
$ws.Range("L6").Value = 11.91773659832167
$ws.Range("N6").Value = 26.10116114395053

$ws.Range("B7").Value = 22.53428645822224
$ws.Range("C7").Value = 3.550255532184509
$ws.Range("D7").Value = 4.602333386737216
$ws.Range("F7").Value = 54.80848217910955
$ws.Range("G7").Value = 3.796990717853251
$ws.Range("J7").Value = 10.86429345800524
$ws.Range("K7").Value = 18.37881109564692
$ws.Range("L7").Value = 11.91260321482487
$ws.Range("N7").Value = 26.08924464976484

$ws.Range("B8").Value = 22.66420126529245
$ws.Range("C8").Value = 3.747615963893701
$ws.Range("D8").Value = 4.591706449981204
$ws.Range("F8").Value = 54.85105527850109
$ws.Range("G8").Value = 3.792787501646844
$ws.Range("J8").Value = 10.83960410869852
$ws.Range("K8").Value = 18.46246447435838
$ws.Range("L8").Value = 11.89388870472836
$ws.Range("N8").Value = 26.04077145431201

$ws.Range("B9").Value = 22.95994904569182
$ws.Range("C9").Value = 4.238621891761761
$ws.Range("D9").Value = 4.575907041787667
$ws.Range("F9").Value = 55.00848951431673
$ws.Range("G9").Value = 3.785361157450185
$ws.Range("J9").Value = 10.79789843242937
$ws.Range("K9").Value = 18.65945539856244
$ws.Range("L9").Value = 11.87132220283266
$ws.Range("N9").Value = 25.9604881716559

$ws.Range("B10").Value = 23.20083443670275
$ws.Range("C10").Value = 4.563237009193221
$ws.Range("D10").Value = 4.567345760829246
$ws.Range("F10").Value = 55.1689164675298
$ws.Range("G10").Value = 3.780396714688746
$ws.Range("J10").Value = 10.77132227179129
$ws.Range("K10").Value = 18.82342240154654
$ws.Range("L10").Value = 11.86333525827713
$ws.Range("N10").Value = 25.91050525141422

$ws.Range("B11").Value = 23.31521508943556
$ws.Range("C11").Value = 4.703007959133656
$ws.Range("D11").Value = 4.56410679549252
$ws.Range("F11").Value = 55.25152310957283
$ws.Range("G11").Value = 3.778243779612408
$ws.Range("J11").Value = 10.76010936578986
$ws.Range("K11").Value = 18.90198744173383
$ws.Range("L11").Value = 11.86156104291459
$ws.Range("N11").Value = 25.88972048342698

$ws.Range("B12").Value = 23.35918854802472
$ws.Range("C12").Value = 4.754797080730023
$ws.Range("D12").Value = 4.562974057233045
$ws.Range("F12").Value = 55.2841780392403
$ws.Range("G12").Value = 3.777443582290193
$ws.Range("J12").Value = 10.75598898351531
$ws.Range("K12").Value = 18.93229069506966
$ws.Range("L12").Value = 11.86115579289233
$ws.Range("N12").Value = 25.88213060495733

$ws.Range("B13").Value = 23.34968921968928
$ws.Range("C13").Value = 4.743694065908228
$ws.Range("D13").Value = 4.563213849362592
$ws.Range("F13").Value = 55.27708431207752
$ws.Range("G13").Value = 3.777615250113846
$ws.Range("J13").Value = 10.75687079679421
$ws.Range("K13").Value = 18.92574009781623
$ws.Range("L13").Value = 11.86123122592679
$ws.Range("N13").Value = 25.88375272857831

$ws.Range("B14").Value = 23.31881976328764
$ws.Range("C14").Value = 4.707291524436993
$ws.Range("D14").Value = 4.564011727629065
$ws.Range("F14").Value = 55.25418218844088
$ws.Range("G14").Value = 3.778177645313428
$ws.Range("J14").Value = 10.75976786231571
$ws.Range("K14").Value = 18.9044695521909
$ws.Range("L14").Value = 11.8615223640601
$ws.Range("N14").Value = 25.88909043028075

$ws.Range("B15").Value = 23.29999635855182
$ws.Range("C15").Value = 4.684845498525585
$ws.Range("D15").Value = 4.564512650329043
$ws.Range("F15").Value = 55.24033251326063
$ws.Range("G15").Value = 3.778524088946594
$ws.Range("J15").Value = 10.76155875776703
$ws.Range("K15").Value = 18.89151207504349
$ws.Range("L15").Value = 11.86173539143083
$ws.Range("N15").Value = 25.89239650521382

$ws.Range("B16").Value = 23.19345342829405
$ws.Range("C16").Value = 4.55394339463594
$ws.Range("D16").Value = 4.567570585899677
$ws.Range("F16").Value = 55.16371076414357
$ws.Range("G16").Value = 3.780539527906631
$ws.Range("J16").Value = 10.77207267179577
$ws.Range("K16").Value = 18.81836636265916
$ws.Range("L16").Value = 11.86348855698017
$ws.Range("N16").Value = 25.91190288384368

$ws.Range("B17").Value = 23.12930116903292
$ws.Range("C17").Value = 4.471613388014344
$ws.Range("D17").Value = 4.569614103690963
$ws.Range("F17").Value = 55.11916436796601
$ws.Range("G17").Value = 3.781802872298661
$ws.Range("J17").Value = 10.77874690319738
$ws.Range("K17").Value = 18.77449866445986
$ws.Range("L17").Value = 11.86503973641704
$ws.Range("N17").Value = 25.92436959609783

$ws.Range("B18").Value = 23.09285624897632
$ws.Range("C18").Value = 4.423516027538598
$ws.Range("D18").Value = 4.570851219661264
$ws.Range("F18").Value = 55.09444939067509
$ws.Range("G18").Value = 3.782539441838036
$ws.Range("J18").Value = 10.78266828781812
$ws.Range("K18").Value = 18.74964255085804
$ws.Range("L18").Value = 11.86610696934323
$ws.Range("N18").Value = 25.93172391375723

$ws.Range("B19").Value = 23.08059546356251
$ws.Range("C19").Value = 4.407103551212245
$ws.Range("D19").Value = 4.571280704168236
$ws.Range("F19").Value = 55.08623739912406
$ws.Range("G19").Value = 3.782790539272832
$ws.Range("J19").Value = 10.78401018946001
$ws.Range("K19").Value = 18.74129174994656
$ws.Range("L19").Value = 11.86649839918318
$ws.Range("N19").Value = 25.93424552431938

$ws.Range("B20").Value = 23.13608353983824
$ws.Range("C20").Value = 4.48045451770729
$ws.Range("D20").Value = 4.569390181504096
$ws.Range("F20").Value = 55.12381261859925
$ws.Range("G20").Value = 3.781667360325886
$ws.Range("J20").Value = 10.77802787966159
$ws.Range("K20").Value = 18.77912973752252
$ws.Range("L20").Value = 11.86485650019305
$ws.Range("N20").Value = 25.92302347092354

$ws.Range("B21").Value = 23.32786920683286
$ws.Range("C21").Value = 4.718014767357812
$ws.Range("D21").Value = 4.563774830024793
$ws.Range("F21").Value = 55.26087191160514
$ws.Range("G21").Value = 3.778012047831678
$ws.Range("J21").Value = 10.75891351536596
$ws.Range("K21").Value = 18.91070239435152
$ws.Range("L21").Value = 11.86142962037213
$ws.Range("N21").Value = 25.88751499484659

$ws.Range("B22").Value = 23.45704425116999
$ws.Range("C22").Value = 4.86663403320794
$ws.Range("D22").Value = 4.560651303624851
$ws.Range("F22").Value = 55.35844843421886
$ws.Range("G22").Value = 3.775710904108607
$ws.Range("J22").Value = 10.74715369050144
$ws.Range("K22").Value = 18.99990279516884
$ws.Range("L22").Value = 11.86074353977478
$ws.Range("N22").Value = 25.86594526553379

$ws.Range("B23").Value = 23.38776079624352
$ws.Range("C23").Value = 4.787921425749707
$ws.Range("D23").Value = 4.562268554517262
$ws.Range("F23").Value = 55.30564190205835
$ws.Range("G23").Value = 3.776931061109868
$ws.Range("J23").Value = 10.75336322642229
$ws.Range("K23").Value = 18.95200776552172
$ws.Range("L23").Value = 11.86096782599065
$ws.Range("N23").Value = 25.87730761767348

$ws.Range("B24").Value = 23.13301586629094
$ws.Range("C24").Value = 4.476459824745641
$ws.Range("D24").Value = 4.569491222735929
$ws.Range("F24").Value = 55.12170835474851
$ws.Range("G24").Value = 3.781728593336195
$ws.Range("J24").Value = 10.7783526876534
$ws.Range("K24").Value = 18.77703489408611
$ws.Range("L24").Value = 11.86493879477199
$ws.Range("N24").Value = 25.92363147147919

$ws.Range("B25").Value = 22.87569015685051
$ws.Range("C25").Value = 4.112108177394001
$ws.Range("D25").Value = 4.57964426881657
$ws.Range("F25").Value = 54.9580105439933
$ws.Range("G25").Value = 3.787283412974433
$ws.Range("J25").Value = 10.80846522532953
$ws.Range("K25").Value = 18.6027196151009
$ws.Range("L25").Value = 11.87591606563936
$ws.Range("N25").Value = 25.98062591561596
